# Adds "The Power User" section to the end of the personas document,
# mirroring the "Gwyn: The Power User" section's Characteristics /
# Interests bullet content, but as a new top-level (Heading1) section.

$d = $word.ActiveDocument

# Paragraphs we will copy the bullet-list numbering format from, so the
# new bullets get their own fresh numId (while keeping the existing
# abstractNumId/bullet look) -- exactly like Word does when a NEW list
# is started from the "Compact" bullet style used elsewhere.
$charsSource = $d.Paragraphs(11)     # "Prefers doing work on the command line..." (numId 1002)
$interestsSource = $d.Paragraphs(15) # "Works in IT, help-desk, linux support" (numId 1003)

function Add-Para($StyleName, $Text) {
    $last = $d.Paragraphs.Last
    $rng = $last.Range
    $rng.Collapse(0) | Out-Null
    $rng.InsertParagraphAfter() | Out-Null
    $p = $d.Paragraphs.Last
    $p.Style = $StyleName
    $r = $p.Range
    $r.Text = $Text
    return $d.Paragraphs.Last
}

function Add-Bookmark($Name) {
    $p = $d.Paragraphs.Last
    $rng = $p.Range
    $d.Bookmarks.Add($Name, $rng) | Out-Null
}

# --- "The Power User" (Heading1) ---
Add-Para "Heading1" "The Power User" | Out-Null
Add-Bookmark "the-power-user"

# --- "Characteristics" (Heading2) ---
Add-Para "Heading2" "Characteristics" | Out-Null
Add-Bookmark "characteristics-3"

# --- Characteristics bullets (new numId, same bullet look as 1002) ---
$p1 = Add-Para "Compact" "Prefers doing work on the command line wherever possible."
$p1.Range.ListFormat.ApplyListTemplate($charsSource.Range.ListFormat.ListTemplate)

$p2 = Add-Para "Compact" "Prefers using the keyboard to control applications."
$p2.Range.ListFormat.ApplyListTemplate($p1.Range.ListFormat.ListTemplate, $true)

$p3 = Add-Para "Compact" "Some understanding of the mechanics of a linux system."
$p3.Range.ListFormat.ApplyListTemplate($p1.Range.ListFormat.ListTemplate, $true)

# --- "Interests" (Heading2) ---
Add-Para "Heading2" "Interests" | Out-Null
Add-Bookmark "interests-1"

# --- Interests bullet (new numId, same bullet look as 1003) ---
$p4 = Add-Para "Compact" "Programs casually."
$p4.Range.ListFormat.ApplyListTemplate($interestsSource.Range.ListFormat.ListTemplate)

# Append a trailing run containing just a space, as its own run (matches
# the source edit's run split) by toggling a direct-character-format on
# just that trailing character so it doesn't coalesce with the run before.
$p4b = $d.Paragraphs.Last
$pEnd = $p4b.Range.End
$p4b.Range.InsertAfter(" ") | Out-Null
$spaceRange = $d.Range($pEnd - 1, $pEnd)
$spaceRange.Font.Bold = $true
$spaceRange.Font.Bold = $false

Write-Output "Added 'The Power User' section."
